{"js": "const replacements = [\n  [\"72\u00f75=\", \"50\u00f76=\"],\n  [\"53\u00f73=\", \"16\u00f76=\"],\n  [\"46\u00f78=\", \"41\u00f74=\"],\n  [\"74\u00f73=\", \"64\u00f77=\"],\n  [\"77\u00f77=\", \"84\u00f76=\"],\n  [\"82\u00f74=\", \"17\u00f72=\"],\n  [\"33\u00f78=\", \"33\u00f73=\"],\n  [\"71\u00f76=\", \"46\u00f79=\"],\n  [\"28\u00f75=\", \"89\u00f72=\"],\n  [\"10\u00f78=\", \"24\u00f77=\"],\n  [\"81\u00f72=\", \"67\u00f72=\"],\n  [\"28\u00f77=\", \"83\u00f74=\"],\n  [\"36\u00f72=\", \"31\u00f74=\"],\n  [\"92\u00f77=\", \"48\u00f73=\"],\n  [\"80\u00f73=\", \"78\u00f76=\"],\n  [\"38\u00f77=\", \"26\u00f79=\"],\n  [\"20\u00f72=\", \"74\u00f78=\"],\n  [\"41\u00f77=\", \"47\u00f76=\"],\n  [\"42\u00f73=\", \"46\u00f74=\"],\n  [\"15\u00f72=\", \"46\u00f76=\"],\n  [\"91\u00f73=\", \"15\u00f76=\"],\n  [\"12\u00f76=\", \"80\u00f74=\"],\n  [\"92\u00f75=\", \"27\u00f79=\"],\n  [\"12\u00f77=\", \"99\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"72\u00f75=\"; New = \"50\u00f76=\" },\n    @{ Old = \"53\u00f73=\"; New = \"16\u00f76=\" },\n    @{ Old = \"46\u00f78=\"; New = \"41\u00f74=\" },\n    @{ Old = \"74\u00f73=\"; New = \"64\u00f77=\" },\n    @{ Old = \"77\u00f77=\"; New = \"84\u00f76=\" },\n    @{ Old = \"82\u00f74=\"; New = \"17\u00f72=\" },\n    @{ Old = \"33\u00f78=\"; New = \"33\u00f73=\" },\n    @{ Old = \"71\u00f76=\"; New = \"46\u00f79=\" },\n    @{ Old = \"28\u00f75=\"; New = \"89\u00f72=\" },\n    @{ Old = \"10\u00f78=\"; New = \"24\u00f77=\" },\n    @{ Old = \"81\u00f72=\"; New = \"67\u00f72=\" },\n    @{ Old = \"28\u00f77=\"; New = \"83\u00f74=\" },\n    @{ Old = \"36\u00f72=\"; New = \"31\u00f74=\" },\n    @{ Old = \"92\u00f77=\"; New = \"48\u00f73=\" },\n    @{ Old = \"80\u00f73=\"; New = \"78\u00f76=\" },\n    @{ Old = \"38\u00f77=\"; New = \"26\u00f79=\" },\n    @{ Old = \"20\u00f72=\"; New = \"74\u00f78=\" },\n    @{ Old = \"41\u00f77=\"; New = \"47\u00f76=\" },\n    @{ Old = \"42\u00f73=\"; New = \"46\u00f74=\" },\n    @{ Old = \"15\u00f72=\"; New = \"46\u00f76=\" },\n    @{ Old = \"91\u00f73=\"; New = \"15\u00f76=\" },\n    @{ Old = \"12\u00f76=\"; New = \"80\u00f74=\" },\n    @{ Old = \"92\u00f75=\"; New = \"27\u00f79=\" },\n    @{ Old = \"12\u00f77=\"; New = \"99\u00f75=\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $r.New, 2) | Out-Null\n}\n"}
